$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H13").Value = 40000
$ws.Range("I13").Value = 40000
$ws.Range("K13").Value = 40000
$ws.Range("M13").Value = -39831
$ws.Range("H19").Value = 1245.8572
$ws.Range("I19").Value = 1066.3334
$ws.Range("J19").Value = 1569
$ws.Range("K19").Value = 1066.3334
$ws.Range("L19").Value = 1569
$ws.Range("M19").Value = -891.3334
$ws.Range("N19").Value = -1919
$ws.Range("H70").Value = 4777.222
$ws.Range("J70").Value = 4777.222
$ws.Range("L70").Value = 14331.666
$ws.Range("N70").Value = -14871.666
$ws.Range("H73").Value = 4777.222
$ws.Range("J73").Value = 4777.222
$ws.Range("L73").Value = 14331.666
$ws.Range("N73").Value = -16203.666
$ws.Range("H74").Value = 0
$ws.Range("J74").Value = 0
$ws.Range("N74").Value = 0
$ws.Range("N74").ClearContents()
$ws.Range("H76").Value = 1800
$ws.Range("J76").Value = 1800
$ws.Range("L76").Value = 1800
$ws.Range("N76").Value = -2430
$ws.Range("H77").Value = 0
$ws.Range("J77").Value = 0
$ws.Range("N77").Value = 0
$ws.Range("N77").ClearContents()
$ws.Range("H79").Value = 1800
$ws.Range("J79").Value = 1800
$ws.Range("L79").Value = 1800
$ws.Range("N79").Value = -3984
$ws.Range("H137").Value = 2193.8
$ws.Range("I137").Value = 1548.6666
$ws.Range("K137").Value = 4645.9998
$ws.Range("M137").Value = -2095.9998
$ws.Range("H138").Value = 2391.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H3").Value = 3541.8
$ws.Range("J3").Value = 3541.8
$ws.Range("L3").Value = 3541.8
$ws.Range("N3").Value = -3771.8
$ws.Range("H45").Value = 1707
$ws.Range("I45").Value = 1707
$ws.Range("K45").Value = 1707
$ws.Range("M45").Value = -1330
$ws.Range("H88").Value = 2326.7144
$ws.Range("I88").Value = 250
$ws.Range("J88").Value = 2486.4614
$ws.Range("K88").Value = 250
$ws.Range("L88").Value = 2486.4614
$ws.Range("M88").Value = 156
$ws.Range("N88").Value = -3298.4614
$ws.Range("H91").Value = 2326.7144
$ws.Range("I91").Value = 250
$ws.Range("J91").Value = 2486.4614
$ws.Range("K91").Value = 250
$ws.Range("L91").Value = 2486.4614
$ws.Range("M91").Value = 1154
$ws.Range("N91").Value = -5294.4614
$ws.Range("H102").Value = 2425.5557
$ws.Range("I102").Value = 2303.75
$ws.Range("K102").Value = 2303.75
$ws.Range("M102").Value = -681.75

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H40").Value = 30000
$ws.Range("J40").Value = 30000
$ws.Range("L40").Value = 30000
$ws.Range("N40").Value = -30530
$ws.Range("H76").Value = 15054.5
$ws.Range("J76").Value = 15054.5
$ws.Range("L76").Value = 15054.5
$ws.Range("N76").Value = -15684.5
$ws.Range("H79").Value = 15054.5
$ws.Range("J79").Value = 15054.5
$ws.Range("L79").Value = 15054.5
$ws.Range("N79").Value = -17238.5
$ws.Range("H96").Value = 28325
$ws.Range("I96").Value = 19950
$ws.Range("J96").Value = 30000
$ws.Range("K96").Value = 19950
$ws.Range("L96").Value = 30000
$ws.Range("M96").Value = -17204
$ws.Range("N96").Value = -35492

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 3569.2856
$ws.Range("I132").Value = 2997
$ws.Range("K132").Value = 8991
$ws.Range("M132").Value = -6461

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 107129.75
$ws.Range("I4").Value = 212000
$ws.Range("J4").Value = 2259.5
$ws.Range("K4").Value = 636000
$ws.Range("L4").Value = 6778.5
$ws.Range("M4").Value = -635888
$ws.Range("N4").Value = -7002.5
$ws.Range("H37").Value = 99521.664
$ws.Range("J37").Value = 99521.664
$ws.Range("L37").Value = 298564.992
$ws.Range("N37").Value = -298788.992
$ws.Range("H98").Value = 281.33334
$ws.Range("J98").Value = 281.33334
$ws.Range("L98").Value = 844.0000200000001
$ws.Range("N98").Value = -3840.00002
$ws.Range("H107").Value = 530.53845
$ws.Range("J107").Value = 456.57144
$ws.Range("L107").Value = 1369.71432
$ws.Range("N107").Value = -5209.71432
$ws.Range("H128").Value = 324893.75
$ws.Range("I128").Value = 324893.75
$ws.Range("K128").Value = 974681.25
$ws.Range("M128").Value = -969701.25
$ws.Range("H138").Value = 4617.5713
$ws.Range("I138").Value = 4248.3335
$ws.Range("K138").Value = 12745.0005
$ws.Range("M138").Value = -7605.000499999998
$ws.Range("H139").Value = 1674.875
$ws.Range("I139").Value = 1316.5
$ws.Range("J139").Value = 2750
$ws.Range("K139").Value = 3949.5
$ws.Range("L139").Value = 8250
$ws.Range("M139").Value = 1190.5
$ws.Range("N139").Value = -18530

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H3").Value = 45472.285
$ws.Range("I3").Value = 2001.5
$ws.Range("J3").Value = 62860.6
$ws.Range("K3").Value = 2001.5
$ws.Range("L3").Value = 62860.6
$ws.Range("M3").Value = -1885.5
$ws.Range("N3").Value = -63092.6
$ws.Range("H5").Value = 0
$ws.Range("I5").Value = 0
$ws.Range("K5").Value = 0
$ws.Range("M5").ClearContents()
$ws.Range("H113").Value = 2074.3635
$ws.Range("I113").Value = 1794
$ws.Range("K113").Value = 1794
$ws.Range("M113").Value = 376
$ws.Range("J122").Value = 1969
$ws.Range("L122").Value = 5907
$ws.Range("N122").Value = -10807

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H4").Value = 10904.5
$ws.Range("I4").Value = 9
$ws.Range("J4").Value = 21800
$ws.Range("K4").Value = 9
$ws.Range("L4").Value = 21800
$ws.Range("M4").Value = 104
$ws.Range("N4").Value = -22026
$ws.Range("H22").Value = 2214.1428
$ws.Range("I22").Value = 500
$ws.Range("K22").Value = 500
$ws.Range("M22").Value = -205
$ws.Range("H27").Value = 2214.1428
$ws.Range("I27").Value = 500
$ws.Range("K27").Value = 500
$ws.Range("M27").Value = -393
$ws.Range("H28").Value = 10904.5
$ws.Range("I28").Value = 9
$ws.Range("J28").Value = 21800
$ws.Range("K28").Value = 9
$ws.Range("L28").Value = 21800
$ws.Range("M28").Value = 223
$ws.Range("N28").Value = -22264
$ws.Range("H37").Value = 10904.5
$ws.Range("I37").Value = 9
$ws.Range("J37").Value = 21800
$ws.Range("K37").Value = 9
$ws.Range("L37").Value = 21800
$ws.Range("M37").Value = 98
$ws.Range("N37").Value = -22014
$ws.Range("H46").Value = 3596.7273
$ws.Range("I46").Value = 2937.7144
$ws.Range("J46").Value = 4750
$ws.Range("K46").Value = 2937.7144
$ws.Range("L46").Value = 4750
$ws.Range("M46").Value = -2749.7144
$ws.Range("N46").Value = -5126
$ws.Range("H55").Value = 182.71428
$ws.Range("I55").Value = 186.33333
$ws.Range("J55").Value = 181.26666
$ws.Range("K55").Value = 186.33333
$ws.Range("L55").Value = 181.26666
$ws.Range("M55").Value = -13.33332999999999
$ws.Range("N55").Value = -527.26666
$ws.Range("H68").Value = 43764.4
$ws.Range("I68").Value = 3414.5
$ws.Range("K68").Value = 3414.5
$ws.Range("M68").Value = -2665.5
$ws.Range("H71").Value = 43764.4
$ws.Range("I71").Value = 3414.5
$ws.Range("K71").Value = 17072.5
$ws.Range("M71").Value = -13328.5
$ws.Range("H82").Value = 400
$ws.Range("J82").Value = 400
$ws.Range("L82").Value = 400
$ws.Range("N82").Value = -1122
$ws.Range("H85").Value = 400
$ws.Range("J85").Value = 400
$ws.Range("L85").Value = 400
$ws.Range("N85").Value = -2896
$ws.Range("H87").Value = 39997
$ws.Range("J87").Value = 39997
$ws.Range("L87").Value = 39997
$ws.Range("N87").Value = -42243
$ws.Range("H90").Value = 39997
$ws.Range("J90").Value = 39997
$ws.Range("L90").Value = 119991
$ws.Range("N90").Value = -131223
$ws.Range("H122").Value = 3557.6667
$ws.Range("I122").Value = 3288.5715
$ws.Range("J122").Value = 4499.5
$ws.Range("K122").Value = 9865.7145
$ws.Range("L122").Value = 13498.5
$ws.Range("M122").Value = -7415.7145
$ws.Range("N122").Value = -18398.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H9").Value = 700
$ws.Range("J9").Value = 700
$ws.Range("L9").Value = 700
$ws.Range("N9").Value = -980
$ws.Range("H12").Value = 700
$ws.Range("J12").Value = 700
$ws.Range("L12").Value = 700
$ws.Range("N12").Value = -984
$ws.Range("H81").Value = 4063.4285
$ws.Range("I81").Value = 4063.4285
$ws.Range("J81").Value = 0
$ws.Range("K81").Value = 8126.857
$ws.Range("L81").Value = 0
$ws.Range("N81").Value = -7065.857
$ws.Range("N81").ClearContents()
$ws.Range("H84").Value = 4063.4285
$ws.Range("I84").Value = 4063.4285
$ws.Range("J84").Value = 0
$ws.Range("K84").Value = 40634.285
$ws.Range("L84").Value = 0
$ws.Range("N84").Value = -35330.285
$ws.Range("N84").ClearContents()
